$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: add year 2022 header in L3, same style as K3
$ws.Range("K3").Copy($ws.Range("L3"))
$ws.Range("L3").Value = 2022

# Row 5: L5 = 12673.2 (population-of-waste total), style like K5
$ws.Range("K5").Copy($ws.Range("L5"))
$ws.Range("L5").Value = 12673.2
$ws.Range("L5").HorizontalAlignment = -4152

# Row 6: L6 = 7037.6 (resident population), style like I6 (no border) + right align
$ws.Range("I6").Copy($ws.Range("L6"))
$ws.Range("L6").Value = 7037.6
$ws.Range("L6").HorizontalAlignment = -4152

# Row 4: L4 = formula, bold right-aligned numeric style like I4 (no border) + bold
$ws.Range("I4").Copy($ws.Range("L4"))
$ws.Range("L4").Formula = "=L5/L6*1000"
$ws.Range("L4").Font.Bold = $true
$ws.Range("L4").HorizontalAlignment = -4152
